$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the rows that changed (re-pulled data / mean calc)
$ws.Range("F4").Value = -4
$ws.Range("F5").Value = 7
$ws.Range("F9").Value = -4
$ws.Range("F12").Value = -3
$ws.Range("F15").Value = 3
$ws.Range("F21").Value = 4
$ws.Range("F22").Value = 0
$ws.Range("F26").Value = -2
